$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (theta_se) updated standard-error values
$ws.Range("B4").Value = "(1.39)"
$ws.Range("C4").Value = "(1.0)"
$ws.Range("D4").Value = "(1.68)"
$ws.Range("E4").Value = "(1.24)"
$ws.Range("F4").Value = "(1.44)"
$ws.Range("G4").Value = "(1.78)"
$ws.Range("H4").Value = "(1.75)"
$ws.Range("I4").Value = "(1.34)"
$ws.Range("J4").Value = "(1.65)"
$ws.Range("K4").Value = "(1.46)"
$ws.Range("L4").Value = "(2.59)"

# Row 6 (lambda_se) updated standard-error values
$ws.Range("B6").Value = "(1.12)"
$ws.Range("C6").Value = "(0.8)"
$ws.Range("D6").Value = "(1.04)"
$ws.Range("E6").Value = "(1.12)"
$ws.Range("F6").Value = "(0.74)"
$ws.Range("G6").Value = "(1.32)"
$ws.Range("H6").Value = "(1.49)"
$ws.Range("I6").Value = "(1.2)"
$ws.Range("J6").Value = "(1.05)"
$ws.Range("K6").Value = "(1.24)"
$ws.Range("L6").Value = "(2.2)"
